$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Date début" / "Date Fin" values on row 9 (E9, F9).
# These are stored as date serials 42865 / 42866 -> 10/05/2017 and 11/05/2017.
$ws.Range("E9").Value = Get-Date -Year 2017 -Month 5 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("F9").Value = Get-Date -Year 2017 -Month 5 -Day 11 -Hour 0 -Minute 0 -Second 0

# Update the selected cell shown in the saved view from E8 to E10.
$ws.Range("E10").Select()
